$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object 'object[,]' 24,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.029925456385346
$arrBF[0,2] = 1.03844555223394
$arrBF[0,3] = 0.9926147277508489
$arrBF[0,4] = 1.045380731931067
$arrBF[1,0] = 1.02
$arrBF[1,1] = 1.030929570257467
$arrBF[1,2] = 1.039245836731627
$arrBF[1,3] = 0.9936372048519304
$arrBF[1,4] = 1.046369565886894
$arrBF[2,0] = 1.02
$arrBF[2,1] = 1.031579122940596
$arrBF[2,2] = 1.039763104252049
$arrBF[2,3] = 0.9942998659930995
$arrBF[2,4] = 1.047009344305518
$arrBF[3,0] = 1.02
$arrBF[3,1] = 1.031852152568096
$arrBF[3,2] = 1.03998042619295
$arrBF[3,3] = 0.9945786998346017
$arrBF[3,4] = 1.047278291324446
$arrBF[4,0] = 1.02
$arrBF[4,1] = 1.031897992957699
$arrBF[4,2] = 1.040016907425374
$arrBF[4,3] = 0.9946255319796338
$arrBF[4,4] = 1.047323447764215
$arrBF[5,0] = 1.02
$arrBF[5,1] = 1.031582771342337
$arrBF[5,2] = 1.039766008659146
$arrBF[5,3] = 0.9943035907982488
$arrBF[5,4] = 1.047012938051819
$arrBF[6,0] = 1.02
$arrBF[6,1] = 1.03026483774921
$arrBF[6,2] = 1.038716129653985
$arrBF[6,3] = 0.9929600610674301
$arrBF[6,4] = 1.045714925690202
$arrBF[7,0] = 1.02
$arrBF[7,1] = 1.027941127264206
$arrBF[7,2] = 1.036861779371546
$arrBF[7,3] = 0.9906006454969559
$arrBF[7,4] = 1.043427206618803
$arrBF[8,0] = 1.02
$arrBF[8,1] = 1.026391088269272
$arrBF[8,2] = 1.035622681834781
$arrBF[8,3] = 0.989033133672735
$arrBF[8,4] = 1.041901790533382
$arrBF[9,0] = 1.02
$arrBF[9,1] = 1.025719690889587
$arrBF[9,2] = 1.035085468643432
$arrBF[9,3] = 0.988355674866747
$arrBF[9,4] = 1.04124121111883
$arrBF[10,0] = 1.02
$arrBF[10,1] = 1.025470270764367
$arrBF[10,2] = 1.034885822847788
$arrBF[10,3] = 0.9881042295826724
$arrBF[10,4] = 1.04099583329396
$arrBF[11,0] = 1.02
$arrBF[11,1] = 1.025523773734786
$arrBF[11,2] = 1.034928652104443
$arrBF[11,3] = 0.9881581567098651
$arrBF[11,4] = 1.041048468080762
$arrBF[12,0] = 1.02
$arrBF[12,1] = 1.025699074423304
$arrBF[12,2] = 1.035068967914339
$arrBF[12,3] = 0.9883348863814464
$arrBF[12,4] = 1.041220928296902
$arrBF[13,0] = 1.02
$arrBF[13,1] = 1.025807078553304
$arrBF[13,2] = 1.035155407764697
$arrBF[13,3] = 0.9884438009545853
$arrBF[13,4] = 1.041327185519043
$arrBF[14,0] = 1.02
$arrBF[14,1] = 1.026435642048575
$arrBF[14,2] = 1.035658320721248
$arrBF[14,3] = 0.9890781214508737
$arrBF[14,4] = 1.041945629697768
$arrBF[15,0] = 1.02
$arrBF[15,1] = 1.026829864246082
$arrBF[15,2] = 1.035973604255818
$arrBF[15,3] = 0.989476357848556
$arrBF[15,4] = 1.042333546511425
$arrBF[16,0] = 1.02
$arrBF[16,1] = 1.027059786142235
$arrBF[16,2] = 1.036157438562567
$arrBF[16,3] = 0.9897087662937556
$arrBF[16,4] = 1.042559805605105
$arrBF[17,0] = 1.02
$arrBF[17,1] = 1.027138179968634
$arrBF[17,2] = 1.036220110224434
$arrBF[17,3] = 0.9897880325774034
$arrBF[17,4] = 1.042636953040781
$arrBF[18,0] = 1.02
$arrBF[18,1] = 1.026787570133261
$arrBF[18,2] = 1.03593978406026
$arrBF[18,3] = 0.9894336180360679
$arrBF[18,4] = 1.042291927348751
$arrBF[19,0] = 1.02
$arrBF[19,1] = 1.025647453658748
$arrBF[19,2] = 1.035027651185837
$arrBF[19,3] = 0.9882828385668249
$arrBF[19,4] = 1.041170143312804
$arrBF[20,0] = 1.02
$arrBF[20,1] = 1.024930424758865
$arrBF[20,2] = 1.034453573591312
$arrBF[20,3] = 0.9875604150241495
$arrBF[20,4] = 1.040464780036992
$arrBF[21,0] = 1.02
$arrBF[21,1] = 1.025310553636935
$arrBF[21,2] = 1.034757958082823
$arrBF[21,3] = 0.9879432794643023
$arrBF[21,4] = 1.040838711352088
$arrBF[22,0] = 1.02
$arrBF[22,1] = 1.026806681087761
$arrBF[22,2] = 1.035955066152411
$arrBF[22,3] = 0.9894529299347244
$arrBF[22,4] = 1.042310733276017
$arrBF[23,0] = 1.02
$arrBF[23,1] = 1.028542020812752
$arrBF[23,2] = 1.037341681550372
$arrBF[23,3] = 0.9912096547607049
$arrBF[23,4] = 1.044018687040146

$ws.Range("B2:F25").Value = $arrBF

$arrIN = New-Object 'object[,]' 24,6
$arrIN[0,0] = 1.034051701522458
$arrIN[0,1] = 1.035069535042062
$arrIN[0,2] = 1.041233564203401
$arrIN[0,3] = 0.9955398523336033
$arrIN[0,4] = 1.048149133117213
$arrIN[0,5] = 1.036539452505985
$arrIN[1,0] = 1.034254408367186
$arrIN[1,1] = 1.035714607601175
$arrIN[1,2] = 1.041843889275786
$arrIN[1,3] = 0.9963617723202692
$arrIN[1,4] = 1.048948933958593
$arrIN[1,5] = 1.037185441142094
$arrIN[2,0] = 1.034383559984822
$arrIN[2,1] = 1.036131260216288
$arrIN[2,2] = 1.042237589546878
$arrIN[2,3] = 0.9968940712668345
$arrIN[2,4] = 1.049465782655221
$arrIN[2,5] = 1.037602685451712
$arrIN[3,0] = 1.034437372845629
$arrIN[3,1] = 1.036306240371293
$arrIN[3,2] = 1.042402808210649
$arrIN[3,3] = 0.997117960005301
$arrIN[3,4] = 1.049682903381083
$arrIN[3,5] = 1.037777914098607
$arrIN[4,0] = 1.034446379954266
$arrIN[4,1] = 1.036335609704951
$arrIN[4,2] = 1.04243053195454
$arrIN[4,3] = 0.9971555583673453
$arrIN[4,4] = 1.04971934936024
$arrIN[4,5] = 1.037807325140088
$arrIN[5,0] = 1.03438428093136
$arrIN[5,1] = 1.036133599019935
$arrIN[5,2] = 1.042239798358611
$arrIN[5,3] = 0.9968970624462087
$arrIN[5,4] = 1.049468684472203
$arrIN[5,5] = 1.037605027576728
$arrIN[6,0] = 1.034120623811555
$arrIN[6,1] = 1.035287695860554
$arrIN[6,2] = 1.041440078816684
$arrIN[6,3] = 0.995817528259106
$arrIN[6,4] = 1.04841956912492
$arrIN[6,5] = 1.036757923137834
$arrIN[7,0] = 1.033640628557577
$arrIN[7,1] = 1.033791363120888
$arrIN[7,2] = 1.040021548474542
$arrIN[7,3] = 0.9939188001724441
$arrIN[7,4] = 1.046565744739352
$arrIN[7,5] = 1.035259465434046
$arrIN[8,0] = 1.033310306649179
$arrIN[8,1] = 1.032789970615905
$arrIN[8,2] = 1.039069629064398
$arrIN[8,3] = 0.9926553831429383
$arrIN[8,4] = 1.04532643438705
$arrIN[8,5] = 1.034256650836844
$arrIN[9,0] = 1.033164829959557
$arrIN[9,1] = 1.032355450748251
$arrIN[9,2] = 1.038655965833062
$arrIN[9,3] = 0.9921088820399291
$arrIN[9,4] = 1.044788992765616
$arrIN[9,5] = 1.033821513901139
$arrIN[10,0] = 1.033110426443823
$arrIN[10,1] = 1.032193914213077
$arrIN[10,2] = 1.038502091599811
$arrIN[10,3] = 0.9919059725120875
$arrIN[10,4] = 1.044589241423623
$arrIN[10,5] = 1.033659747965556
$arrIN[11,0] = 1.033122112784313
$arrIN[11,1] = 1.03222857052156
$arrIN[11,2] = 1.038535108152121
$arrIN[11,3] = 0.9919494934313052
$arrIN[11,4] = 1.044632094260482
$arrIN[11,5] = 1.033694453489973
$arrIN[12,0] = 1.033160340432544
$arrIN[12,1] = 1.03234210087117
$arrIN[12,2] = 1.038643251052706
$arrIN[12,3] = 0.9920921077337197
$arrIN[12,4] = 1.044772483735529
$arrIN[12,5] = 1.033808145065702
$arrIN[13,0] = 1.033183845128316
$arrIN[13,1] = 1.03241203258519
$arrIN[13,2] = 1.038709852152688
$arrIN[13,3] = 0.9921799884222134
$arrIN[13,4] = 1.044858966224195
$arrIN[13,5] = 1.033878176090777
$arrIN[14,0] = 1.033319909972586
$arrIN[14,1] = 1.032818789111211
$arrIN[14,2] = 1.039097051473969
$arrIN[14,3] = 0.9926916645766087
$arrIN[14,4] = 1.045362085543363
$arrIN[14,5] = 1.034285510257721
$arrIN[15,0] = 1.033404605479501
$arrIN[15,1] = 1.033073693200815
$arrIN[15,2] = 1.039339536565607
$arrIN[15,3] = 0.9930127773699352
$arrIN[15,4] = 1.0456774618801
$arrIN[15,5] = 1.034540776340369
$arrIN[16,0] = 1.033453770993296
$arrIN[16,1] = 1.033222286595855
$arrIN[16,2] = 1.039480831590566
$arrIN[16,3] = 0.9932001317071769
$arrIN[16,4] = 1.045861337182633
$arrIN[16,5] = 1.034689580755075
$arrIN[17,0] = 1.03347049513207
$arrIN[17,1] = 1.033272938195225
$arrIN[17,2] = 1.039528985371204
$arrIN[17,3] = 0.9932640239640975
$arrIN[17,4] = 1.045924020602024
$arrIN[17,5] = 1.034740304285525
$arrIN[18,0] = 1.033395542855148
$arrIN[18,1] = 1.033046353488951
$arrIN[18,2] = 1.039313534938683
$arrIN[18,3] = 0.9929783193494215
$arrIN[18,4] = 1.045643633087875
$arrIN[18,5] = 1.034513397802979
$arrIN[19,0] = 1.033149093470994
$arrIN[19,1] = 1.032308672779205
$arrIN[19,2] = 1.038611411769112
$arrIN[19,3] = 0.9920501090198102
$arrIN[19,4] = 1.044731145887603
$arrIN[19,5] = 1.033774669502011
$arrIN[20,0] = 1.032992017808138
$arrIN[20,1] = 1.031844073957932
$arrIN[20,2] = 1.038168678954989
$arrIN[20,3] = 0.9914670000341481
$arrIN[20,4] = 1.044156725657341
$arrIN[20,5] = 1.033309410897122
$arrIN[21,0] = 1.033075487722299
$arrIN[21,1] = 1.032090441353472
$arrIN[21,2] = 1.038403501192468
$arrIN[21,3] = 0.991776070289318
$arrIN[21,4] = 1.044461303221583
$arrIN[21,5] = 1.033556128162622
$arrIN[22,0] = 1.03339963859399
$arrIN[22,1] = 1.033058707399268
$arrIN[22,2] = 1.039325284394671
$arrIN[22,3] = 0.9929938892766442
$arrIN[22,4] = 1.045658919105032
$arrIN[22,5] = 1.034525769257265
$arrIN[23,0] = 1.033766540683544
$arrIN[23,1] = 1.034178879093428
$arrIN[23,2] = 1.04038937363632
$arrIN[23,3] = 0.9944092447426414
$arrIN[23,4] = 1.047045609375302
$arrIN[23,5] = 1.035647531723716

$ws.Range("I2:N25").Value = $arrIN

Write-Host "Update complete"
